# ncp-gop-transect-summer-2018/ncplterEn617.xlsx
# Remove the "temp" (row 7) and "sal" (row 8) attribute rows from the
# ColumnHeaders metadata sheet, shifting the rows below (biosat,
# O2_Ar_ratio_corrected, ncp, k) up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows in one shot so everything below shifts up and the
# used-range dimension shrinks from G12 to G10, matching the authored diff.
$ws.Rows("7:8").Delete()

# Restore a plain view on the sheet: no frozen/scrolled top-left cell, and
# the active selection parked on A15 (an empty row below the shrunk table),
# as in the post-edit workbook.
$ws.Range("A15").Select()
